$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link columns (B, C) -- plain text updates.
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"

# Price / volume columns (D, E) hold numeric-looking text (e.g. "331.99",
# "0.60%") that must stay plain text, matching the original inlineStr cells,
# instead of being auto-converted to numbers/percentages by Excel. Force the
# cell format to Text ("@") immediately before writing each value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.93%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.770"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.25%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08063"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.76%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.506"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.634"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.44%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.959"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.66%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.73%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9222"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.87%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1287"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.09%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1954"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.15%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.731"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "17.24%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09310"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.88%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1052"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.75%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001310"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.79%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006270"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.24%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.370"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.12%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3484"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.96%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1341"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.01%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2671"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.92%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04430"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.63%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.87%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004509"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.96%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001202"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02518"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.10%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05469"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.03%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007494"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.74%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009918"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "12.39%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1410"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.49%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002110"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01132"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "17.86%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006801"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.76%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003033"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.73%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002282"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.91%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.13%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.13%"
